$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "28.351.12"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +1.60%  "

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.825.97"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +2.85%  "

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.001"
$r.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "317.38"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "1.002"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.5348"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.4063"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +9.18%  "

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.07599"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +2.62%  "

$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("E11").Value = "  +1.87%  "

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "6.332"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +4.62%  "

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "1.003"
$r.Style = "Normal"
$ws.Range("E13").Value = "  +0.05%  "

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "7.608"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +5.76%  "

$ws.Range("E15").Value = "  +1.88%  "

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "1.824.38"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +3.46%  "

$ws.Range("E17").Value = "  +1.91%  "

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.00001073"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +2.22%  "

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.06611"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +3.30%  "

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "17.62"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +2.15%  "

$ws.Range("E21").Value = "  -0.02%  "

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "6.091"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +3.79%  "

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "28.362.97"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +1.58%  "

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "11.19"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +1.19%  "

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "2.182"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +5.18%  "

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "2.464"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +8.74%  "

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "157.94"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +1.49%  "

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "20.56"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +2.09%  "

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "2.037.66"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +3.41%  "

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "123.91"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +3.58%  "

$ws.Range("E31").Value = "  +1.62%  "

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "0.1095"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +5.62%  "

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "5.649"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +3.03%  "

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "3.647"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +0.39%  "

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.07275"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +15.06%  "

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.2244"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +0.79%  "

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.02343"
$r.Style = "Normal"
$ws.Range("E37").Value = "  +3.88%  "

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "5.190"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +5.02%  "

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "8.860"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +6.04%  "

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.6244"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +2.37%  "

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "11.27"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +3.23%  "

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "1.184"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("E43").Value = "  -0.04%  "

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "1.400"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -1.87%  "

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "13.47"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +2.23%  "

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "3.705"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +1.40%  "

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.5841"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +2.09%  "

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "125.47"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("E49").Value = "  +3.53%  "

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "1.204"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +2.23%  "

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.06892"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +1.49%  "
